# "added a working version of buy from pega"
# Update the "Shopping List" sheet: change row 5's item/quantity/order-id,
# and append three new purchase rows (6-8) pulled from the Products list.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Shopping List")

# Row 5: was "beer" / 3 / 12345678 -> now "Outback Lager" / 5 / 47414620
$ws.Range("A5").Value = "Outback Lager"
$ws.Range("B5").Value = 5
$ws.Range("C5").Value = 127508
$ws.Range("D5").Value = 47414620

# Row 6 (new): "Chai" (reuses the now-freed "beer" shared string slot) / 3 / 37614806
$ws.Range("A6").Value = "Chai"
$ws.Range("B6").Value = 3
$ws.Range("C6").Value = 127508
$ws.Range("D6").Value = 37614806

# Row 7 (new): "Veggie Spread" / 2 / 73752320
$ws.Range("A7").Value = "Veggie Spread"
$ws.Range("B7").Value = 2
$ws.Range("C7").Value = 127508
$ws.Range("D7").Value = 73752320

# Row 8 (new): "Aniseed Syrup" / 2 / 36295933
$ws.Range("A8").Value = "Aniseed Syrup"
$ws.Range("B8").Value = 2
$ws.Range("C8").Value = 127508
$ws.Range("D8").Value = 36295933

# Move the active selection up from A6 to A5, matching the saved view state
$ws.Range("A5").Select() | Out-Null
